$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new draw-result row (row 50) exactly as the upstream
# auto-updater writes it: every column is literal text, including the
# date-looking, all-digit, and timestamp-looking values. Prefixing the
# ambiguous ones with an apostrophe forces Excel to store them as text
# instead of auto-converting to a date/number serial.
$ws.Range("A50").Value = "'2025-11-05"
$ws.Range("B50").Value = "Pick 4"
$ws.Range("C50").Value = "'251105"
$ws.Range("D50").Value = "7-5-9-7"
$ws.Range("E50").Value = "'2025-11-05T21:39:17.374+04:00"
